$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows down), to hold the newest date entry
$ws.Rows("2:2").Insert(1)

# Populate the new row 2 with the latest date and the same price values as the rest of the dataset.
# Force the date cell to be stored as literal text (matching the rest of the column) instead of
# letting Excel auto-convert the "yyyy-mm-dd" looking string into a date serial number.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2026-02-14"
$ws.Range("A2").NumberFormat = "General"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 783.5
$ws.Range("C2").Value = 1112
$ws.Range("D2").Value = 3610
